# Applies the "DRI EAR valuesV2" edit:
#  - renames the nutrient code for Vitamin A from "vitamin_a_IU" to "vitamin_a"
#  - adds two summary columns (Y = min, Z = max) across the DRI table (rows 1, 3-27)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Add header labels for the new "min"/"max" columns
# ---------------------------------------------------------------------------
$ws.Range("Y1").Value = "min"
$ws.Range("Z1").Value = "max"

# ---------------------------------------------------------------------------
# 2) Rename the "vitamin_a_IU" nutrient code to "vitamin_a" (row 8, column A)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "vitamin_a"

# ---------------------------------------------------------------------------
# 3) Add MIN / MAX formulas for each data row (row 3 through row 27), spanning
#    the age/sex bracket columns C:X
# ---------------------------------------------------------------------------
for ($r = 3; $r -le 27; $r++) {
    $minFormula = "=MIN(C" + $r + ":X" + $r + ")"
    $maxFormula = "=MAX(C" + $r + ":X" + $r + ")"
    $ws.Range("Y$r").Formula = $minFormula
    $ws.Range("Z$r").Formula = $maxFormula
}

# Match the existing "0.0" numeric style used throughout the table (style 13)
$ws.Range("Y3:Z27").NumberFormat = "0.0"
